# The deck ships two DrawingML theme parts:
#   theme1.xml -> currently "Office Theme" colours (used by the Notes Master)
#   theme2.xml -> currently "Integral" colours      (used by the Slide Master /
#                                                     every slide's live theme)
#
# The authored edit swaps the two themes' contents, so the slides (and the
# presentation's active design) end up using the stock "Office Theme" palette
# instead of "Integral". Re-create that visual result by rewriting each of the
# twelve theme colour slots exposed on the live ThemeColorScheme (dk1, lt1,
# dk2, lt2, accent1-6, hlink, folHlink) to the Office Theme's RGB values.

function Convert-RGB($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

$tcs.Colors(1).RGB  = Convert-RGB 0x00 0x00 0x00   # dk1
$tcs.Colors(2).RGB  = Convert-RGB 0xFF 0xFF 0xFF   # lt1
$tcs.Colors(3).RGB  = Convert-RGB 0x44 0x54 0x6A   # dk2
$tcs.Colors(4).RGB  = Convert-RGB 0xE7 0xE6 0xE6   # lt2
$tcs.Colors(5).RGB  = Convert-RGB 0x5B 0x9B 0xD5   # accent1
$tcs.Colors(6).RGB  = Convert-RGB 0xED 0x7D 0x31   # accent2
$tcs.Colors(7).RGB  = Convert-RGB 0xA5 0xA5 0xA5   # accent3
$tcs.Colors(8).RGB  = Convert-RGB 0xFF 0xC0 0x00   # accent4
$tcs.Colors(9).RGB  = Convert-RGB 0x44 0x72 0xC4   # accent5
$tcs.Colors(10).RGB = Convert-RGB 0x70 0xAD 0x47   # accent6
$tcs.Colors(11).RGB = Convert-RGB 0x05 0x63 0xC1   # hlink
$tcs.Colors(12).RGB = Convert-RGB 0x95 0x4F 0x72   # folHlink
